$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Cells.Item(2, 4) "57.691.64"
$ws.Cells.Item(2, 5).Value = "  +0.34%  "
Set-TextValue $ws.Cells.Item(3, 4) "3.122.33"
$ws.Cells.Item(3, 5).Value = "  +0.44%  "
$ws.Cells.Item(4, 5).Value = "  +0.01%  "
Set-TextValue $ws.Cells.Item(5, 4) "532.59"
$ws.Cells.Item(5, 5).Value = "  +1.20%  "
Set-TextValue $ws.Cells.Item(6, 4) "138.07"
$ws.Cells.Item(6, 5).Value = "  +0.75%  "
$ws.Cells.Item(7, 5).Value = "  -0.01%  "
Set-TextValue $ws.Cells.Item(8, 4) "3.120.16"
$ws.Cells.Item(8, 5).Value = "  +0.45%  "
$ws.Cells.Item(9, 5).Value = "  +6.46%  "
$ws.Cells.Item(10, 5).Value = "  +0.35%  "
$ws.Cells.Item(11, 5).Value = "  +0.12%  "
$ws.Cells.Item(12, 5).Value = "  +4.55%  "
$ws.Cells.Item(13, 5).Value = "  +1.51%  "
Set-TextValue $ws.Cells.Item(14, 4) "3.660.37"
$ws.Cells.Item(14, 5).Value = "  +0.43%  "
Set-TextValue $ws.Cells.Item(15, 4) "25.95"
$ws.Cells.Item(15, 5).Value = "  +2.40%  "
$ws.Cells.Item(16, 5).Value = "  +0.54%  "
Set-TextValue $ws.Cells.Item(17, 4) "57.798.87"
$ws.Cells.Item(17, 5).Value = "  +0.35%  "
Set-TextValue $ws.Cells.Item(18, 4) "3.123.32"
$ws.Cells.Item(18, 5).Value = "  +0.53%  "
Set-TextValue $ws.Cells.Item(19, 4) "6.07"
$ws.Cells.Item(19, 5).Value = "  +2.19%  "
$ws.Cells.Item(20, 5).Value = "  +2.42%  "
$ws.Cells.Item(21, 5).Value = "  +2.42%  "
Set-TextValue $ws.Cells.Item(22, 4) "368.22"
$ws.Cells.Item(22, 5).Value = "  +6.58%  "
Set-TextValue $ws.Cells.Item(23, 4) "0.999"
$ws.Cells.Item(23, 5).Value = "  -0.15%  "
$ws.Cells.Item(24, 5).Value = "  -2.08%  "
Set-TextValue $ws.Cells.Item(25, 4) "69.12"
$ws.Cells.Item(25, 5).Value = "  +2.09%  "
$ws.Cells.Item(26, 5).Value = "  +1.08%  "
$ws.Cells.Item(27, 5).Value = "  +0.21%  "
$ws.Cells.Item(28, 5).Value = "  +0.08%  "
Set-TextValue $ws.Cells.Item(29, 4) "0.0₃0862"
$ws.Cells.Item(29, 5).Value = "  -3.60%  "
$ws.Cells.Item(30, 5).Value = "  -1.81%  "
$ws.Cells.Item(31, 5).Value = "  -0.04%  "
$ws.Cells.Item(32, 5).Value = "  +0.41%  "
$ws.Cells.Item(33, 5).Value = "  +2.57%  "
Set-TextValue $ws.Cells.Item(34, 4) "5.14"
$ws.Cells.Item(34, 5).Value = "  +3.67%  "
$ws.Cells.Item(35, 5).Value = "  -0.62%  "
Set-TextValue $ws.Cells.Item(36, 4) "159.44"
$ws.Cells.Item(36, 5).Value = "  +0.62%  "
Set-TextValue $ws.Cells.Item(37, 4) "6.07"
$ws.Cells.Item(37, 5).Value = "  +0.02%  "
$ws.Cells.Item(38, 5).Value = "  +5.08%  "
Set-TextValue $ws.Cells.Item(39, 4) "25.38"
$ws.Cells.Item(39, 5).Value = "  -1.93%  "
$ws.Cells.Item(40, 5).Value = "  +3.57%  "
$ws.Cells.Item(41, 5).Value = "  +1.56%  "
Set-TextValue $ws.Cells.Item(42, 4) "2.539.54"
$ws.Cells.Item(42, 5).Value = "  +6.76%  "
Set-TextValue $ws.Cells.Item(43, 4) "4.06"
$ws.Cells.Item(43, 5).Value = "  -0.99%  "
$ws.Cells.Item(44, 5).Value = "  -0.25%  "
Set-TextValue $ws.Cells.Item(45, 4) "37.80"
$ws.Cells.Item(45, 5).Value = "  +3.26%  "
Set-TextValue $ws.Cells.Item(46, 4) "0.999"
$ws.Cells.Item(46, 5).Value = "  -0.05%  "
$ws.Cells.Item(47, 5).Value = "  +0.57%  "
Set-TextValue $ws.Cells.Item(48, 4) "0.977"
$ws.Cells.Item(48, 5).Value = "  +0.25%  "
$ws.Cells.Item(49, 5).Value = "  +2.00%  "
Set-TextValue $ws.Cells.Item(50, 4) "19.67"
$ws.Cells.Item(50, 5).Value = "  -0.63%  "
Set-TextValue $ws.Cells.Item(51, 4) "0.739"
$ws.Cells.Item(51, 5).Value = "  -2.99%  "
